$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.090.59'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.617.07'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.97'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.98%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '658.81'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.66%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.73'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +20.05%  '
$ws.Range("E8").Value = '  +5.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.07'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +9.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.611.42'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.01'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.95%  '
$ws.Range("E13").Value = '  +2.47%  '
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.287.83'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '97.942.07'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.92%  '
$ws.Range("E17").Value = '  +3.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.612.98'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.03'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.18%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.90'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.10'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.532'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +13.53%  '
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '514.63'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.47%  '
$ws.Range("E25").Value = '  +8.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.90'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '99.69'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +5.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.808.03'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.155'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +12.53%  '
$ws.Range("E31").Value = '  +1.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.81'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +5.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.187'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +5.47%  '
$ws.Range("E35").Value = '  -0.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.91'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.85'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +8.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.572'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '611.62'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +9.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.63'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +9.13%  '
$ws.Range("E41").Value = '  +13.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.153'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.99%  '
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("E44").Value = '  +3.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.99'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +7.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0441'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +8.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.31'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.23%  '
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.62'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +7.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.400'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +37.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.03'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.82%  '
